$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the weekly figures for each team (rows 2-8) ---
# Columns: B = Team, C = Weekly Pending Total(Rp), D = Repayment
# Column E holds the (shared) formula =D/C which recalculates automatically.

$data = @(
    @{ Row = 2;  Team = "Kurni_s2l";    C = 6832283586; D = 724401376 },
    @{ Row = 3;  Team = "Xinghao_s2l";  C = 2817349113; D = 294469421 },
    @{ Row = 4;  Team = "Cpu_s2l";      C = 1872312028; D = 192509595 },
    @{ Row = 5;  Team = "EDN_S2l";      C = 1866206765; D = 182704064 },
    @{ Row = 6;  Team = "Mkm_s2l";      C = 3663615678; D = 336547776 },
    @{ Row = 7;  Team = "Hansyah_S2l";  C = 6908517958; D = 627734875 },
    @{ Row = 8;  Team = "Zakka_S2l";    C = 6770867666; D = 606553043 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.Team
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
}

# --- Update the saved selection to the recovery-rate column ---
$ws.Range("E2:E8").Select()
